# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 70 (pushing the existing rows 70-110
# down to 71-111) and populate it with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 70; Excel shifts rows 70..110
# down to 71..111 automatically (and the sheet dimension grows to R111).
$ws.Rows("70:70").Insert()

$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 44529
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = 100112045
$ws.Range("G70").Value = "Zapallo"
$ws.Range("H70").Value = "Paine"
$ws.Range("I70").Value = "1a (guarda)"
$ws.Range("J70").Value = 160
$ws.Range("K70").Value = 220
$ws.Range("L70").Value = 250
$ws.Range("M70").Value = 235
$ws.Range("N70").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O70").Value = "Región de O'Higgins"
$ws.Range("P70").Value = 235
$ws.Range("Q70").Value = 1
$ws.Range("R70").Value = "Hortaliza"
